$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell Y1 - week 22 label, styled like the other week-number headers (bold, centered)
$ws.Range("Y1").Font.Bold = $true
$ws.Range("Y1").HorizontalAlignment = -4108
$ws.Range("Y1").Value = "'22"

# Weekly counts for week 22 (new column Y)
$ws.Range("Y2").Value = 43
$ws.Range("Y5").Value = 1
$ws.Range("Y6").Value = 35
$ws.Range("Y7").Value = 32
$ws.Range("Y8").Value = 27
$ws.Range("Y10").Value = 2
$ws.Range("Y11").Value = 2
$ws.Range("Y12").Value = 2
$ws.Range("Y13").Value = 4
$ws.Range("Y14").Value = 1
$ws.Range("Y15").Value = 2
$ws.Range("Y16").Value = 4
$ws.Range("Y17").Value = 1
$ws.Range("Y19").Value = 1
$ws.Range("Y21").Value = 2
$ws.Range("Y22").Value = 6
$ws.Range("Y23").Value = 1
$ws.Range("Y24").Value = 51
$ws.Range("Y27").Value = 282
$ws.Range("Y30").Value = 3
$ws.Range("Y31").Value = 8
$ws.Range("Y33").Value = 0
$ws.Range("Y34").Value = 26
$ws.Range("Y36").Value = 14
$ws.Range("Y37").Value = 99
$ws.Range("Y38").Value = 3
$ws.Range("Y39").Value = 9
$ws.Range("Y40").Value = 74
$ws.Range("Y41").Value = 41
$ws.Range("Y42").Value = 186
$ws.Range("Y43").Value = 73
$ws.Range("Y44").Value = 98
$ws.Range("Y45").Value = 0
$ws.Range("Y46").Value = 66
$ws.Range("Y47").Value = 3
$ws.Range("Y48").Value = 0
$ws.Range("Y49").Value = 14
$ws.Range("Y51").Value = 49
$ws.Range("Y52").Value = 0
$ws.Range("Y53").Value = 0
$ws.Range("Y54").Value = 2
$ws.Range("Y55").Value = 19
$ws.Range("Y56").Value = 31
